# Append a new data row (row 32) to Sheet1, columns A and B.
# The source values are numeric-looking but must be stored as literal TEXT
# (e.g. "35600.0" keeps its trailing ".0"), so each value is entered with a
# leading apostrophe - the standard Excel "force text" quote prefix - which
# prevents auto-conversion to a Number. The Style is reset to "Normal"
# afterwards so the cells keep the default (General) display style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellA32 = $ws.Cells.Item(32, 1)
$cellA32.Value = "'35600.0"
$cellA32.Style = "Normal"

$cellB32 = $ws.Cells.Item(32, 2)
$cellB32.Value = "'64641.0"
$cellB32.Style = "Normal"
